$d = $word.ActiveDocument

# 1. Generalize the plaintiff merge-field: "{{ plaintiff_name }}" -> "{{ plaintiffs }}"
#    (drop the spell-check markup around the old "plaintiff_name" token as part of the rewrite)
$d.Content.Find.Execute("{{ plaintiff_name}}", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "{{ plaintiffs }}", 2) | Out-Null

# 2. Generalize the defendant merge-field: "{{ defendant_name }}" -> "{{ defendants }}"
$d.Content.Find.Execute("{{ defendant_name }}", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "{{ defendants }}", 2) | Out-Null

# 3. Merge the (previously split) double-underlined title runs
#    " " + "PURSUANT TO MASS. R. A. P. 6(A)" -> " PURSUANT TO MASS. R. A. P. 6(A)"
$titleRange = $d.Content
$found = $titleRange.Find.Execute("PURSUANT TO MASS. R. A. P. 6(A)")
if ($found) {
    $mergedRange = $d.Range($titleRange.Start - 1, $titleRange.End)
    # touch the range with different text first so the two like-formatted runs
    # are forced to coalesce into a single run once the real text is restored
    $mergedRange.Text = " PURSUANT TO MASS. R. A. P. 6(A)#"
    $restoreRange = $d.Range($titleRange.Start - 1, $titleRange.End + 1)
    $restoreRange.Text = " PURSUANT TO MASS. R. A. P. 6(A)"
}
